$wb = $excel.ActiveWorkbook

# ALC row 2 (@@ -727,22 +727,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1846.0667
$ws.Range("I2").Value = 1636.5
$ws.Range("K2").Value = 1636.5
$ws.Range("M2").Value = -1523.5

# ALC row 86 (@@ -4909,25 +4909,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 3108.75
$ws.Range("I86").Value = 2767.5386
$ws.Range("J86").Value = 3742.4285
$ws.Range("K86").Value = 2767.5386
$ws.Range("L86").Value = 3742.4285
$ws.Range("M86").Value = -1644.5386
$ws.Range("N86").Value = -5988.4285

# ALC row 89 (@@ -5062,25 +5062,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 3108.75
$ws.Range("I89").Value = 2767.5386
$ws.Range("J89").Value = 3742.4285
$ws.Range("K89").Value = 13837.693
$ws.Range("L89").Value = 18712.1425
$ws.Range("M89").Value = -8221.692999999999
$ws.Range("N89").Value = -29944.1425

# ALC row 112 (@@ -6222,25 +6222,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 6767.9346
$ws.Range("J112").Value = 7181.3257
$ws.Range("L112").Value = 21543.9771
$ws.Range("N112").Value = -23759.9771

# ALC row 138 (@@ -7529,25 +7529,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2736.5
$ws.Range("I138").Value = 2112.2
$ws.Range("J138").Value = 3985.1
$ws.Range("K138").Value = 6336.599999999999
$ws.Range("L138").Value = 11955.3
$ws.Range("M138").Value = -1196.599999999999
$ws.Range("N138").Value = -22235.3

# ARM row 32 (@@ -9301,22 +9301,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18110.8
$ws.Range("I32").Value = 18342.041
$ws.Range("J32").Value = 999
$ws.Range("K32").Value = 18342.041
$ws.Range("L32").Value = 999
$ws.Range("M32").Value = -18055.041
$ws.Range("N32").Value = -1573

# ARM row 45 (@@ -9932,25 +9935,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3322.6365
$ws.Range("I45").Value = 2103.1765
$ws.Range("J45").Value = 4618.3125
$ws.Range("K45").Value = 2103.1765
$ws.Range("L45").Value = 4618.3125
$ws.Range("M45").Value = -1726.1765
$ws.Range("N45").Value = -5372.3125

# ARM row 132 (@@ -14138,25 +14141,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2562.2144
$ws.Range("I132").Value = 2780
$ws.Range("J132").Value = 2271.8333
$ws.Range("K132").Value = 8340
$ws.Range("L132").Value = 6815.499899999999
$ws.Range("M132").Value = -5810
$ws.Range("N132").Value = -11875.4999

# BSM row 19 (@@ -15564,19 +15567,22 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 1000
$ws.Range("I19").Value = 1000
$ws.Range("K19").Value = 1000
$ws.Range("M19").Value = -827

# BSM row 20 (@@ -15610,22 +15616,22 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 22166.25
$ws.Range("I20").Value = 27549.525
$ws.Range("K20").Value = 27549.525
$ws.Range("M20").Value = -27302.525

# BSM row 23 (@@ -15763,19 +15769,22 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 2000
$ws.Range("J23").Value = 2000
$ws.Range("L23").Value = 2000
$ws.Range("N23").Value = -2566

# BSM row 86 (@@ -18844,25 +18853,25 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1263.3478
$ws.Range("I86").Value = 1286.5555
$ws.Range("J86").Value = 1179.8
$ws.Range("K86").Value = 1286.5555
$ws.Range("L86").Value = 1179.8
$ws.Range("M86").Value = -163.5554999999999
$ws.Range("N86").Value = -3425.8

# BSM row 89 (@@ -18994,25 +19003,25 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1263.3478
$ws.Range("I89").Value = 1286.5555
$ws.Range("J89").Value = 1179.8
$ws.Range("K89").Value = 6432.7775
$ws.Range("L89").Value = 5899
$ws.Range("M89").Value = -816.7775000000001
$ws.Range("N89").Value = -17131

# BSM row 105 (@@ -19790,22 +19799,22 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3193.0908
$ws.Range("I105").Value = 890.625
$ws.Range("K105").Value = 890.625
$ws.Range("M105").Value = 856.375

# CRP row 31 (@@ -23085,22 +23094,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I31").Value = 3336190
$ws.Range("J31").Value = 9004
$ws.Range("K31").Value = 3336190
$ws.Range("L31").Value = 9004
$ws.Range("M31").Value = -3335895
$ws.Range("N31").Value = -9594

# CRP row 34 (@@ -23232,22 +23241,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I34").Value = 3336190
$ws.Range("J34").Value = 9004
$ws.Range("K34").Value = 3336190
$ws.Range("L34").Value = 9004
$ws.Range("M34").Value = -3335988
$ws.Range("N34").Value = -9408

# CRP row 127 (@@ -27822,22 +27831,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H127").Value = 81427.71000000001
$ws.Range("J127").Value = 81427.71000000001
$ws.Range("L127").Value = 81427.71000000001
$ws.Range("N127").Value = -91347.71000000001

# CRP row 132 (@@ -28064,22 +28073,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 48253.23
$ws.Range("I132").Value = 61718.3
$ws.Range("K132").Value = 185154.9
$ws.Range("M132").Value = -182624.9

# CUL row 113 (@@ -34246,25 +34255,25 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2388.2273
$ws.Range("I113").Value = 3174.1428
$ws.Range("J113").Value = 2021.4667
$ws.Range("K113").Value = 9522.428400000001
$ws.Range("L113").Value = 6064.4001
$ws.Range("M113").Value = -7352.428400000001
$ws.Range("N113").Value = -10404.4001

# CUL row 122 (@@ -34702,25 +34711,25 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 1160.2858
$ws.Range("J122").Value = 1253.5
$ws.Range("L122").Value = 11281.5
$ws.Range("N122").Value = -16181.5

# GSM row 80 (@@ -39613,22 +39622,22 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I80").Value = 4687.5
$ws.Range("J80").Value = 14371
$ws.Range("K80").Value = 4687.5
$ws.Range("L80").Value = 14371
$ws.Range("M80").Value = -3689.5
$ws.Range("N80").Value = -16367

# GSM row 83 (@@ -39763,22 +39772,22 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I83").Value = 4687.5
$ws.Range("J83").Value = 14371
$ws.Range("K83").Value = 23437.5
$ws.Range("L83").Value = 71855
$ws.Range("M83").Value = -18445.5
$ws.Range("N83").Value = -81839

# GSM row 107 (@@ -40918,25 +40927,25 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 269.9091
$ws.Range("I107").Value = 165.125
$ws.Range("J107").Value = 549.3333
$ws.Range("K107").Value = 165.125
$ws.Range("L107").Value = 549.3333
$ws.Range("M107").Value = 1754.875
$ws.Range("N107").Value = -4389.3333

# GSM row 113 (@@ -41212,25 +41221,25 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1854.875
$ws.Range("I113").Value = 1894.65
$ws.Range("J113").Value = 1656
$ws.Range("K113").Value = 1894.65
$ws.Range("L113").Value = 1656
$ws.Range("M113").Value = 275.3499999999999
$ws.Range("N113").Value = -5996

# GSM row 132 (@@ -42131,25 +42140,25 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2098.1667
$ws.Range("I132").Value = 1631.5
$ws.Range("J132").Value = 2331.5
$ws.Range("K132").Value = 4894.5
$ws.Range("L132").Value = 6994.5
$ws.Range("M132").Value = -2364.5
$ws.Range("N132").Value = -12054.5

# LTW row 2 (@@ -42718,23 +42727,20 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 4000000
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()

# LTW row 22 (@@ -43725,25 +43731,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1170.7084
$ws.Range("I22").Value = 807.53845
$ws.Range("J22").Value = 1599.909
$ws.Range("K22").Value = 807.53845
$ws.Range("L22").Value = 1599.909
$ws.Range("M22").Value = -512.53845
$ws.Range("N22").Value = -2189.909

# LTW row 27 (@@ -43982,25 +43988,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1170.7084
$ws.Range("I27").Value = 807.53845
$ws.Range("J27").Value = 1599.909
$ws.Range("K27").Value = 807.53845
$ws.Range("L27").Value = 1599.909
$ws.Range("M27").Value = -700.53845
$ws.Range("N27").Value = -1813.909

# WVR row 11 (@@ -50158,22 +50164,22 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 17000
$ws.Range("J11").Value = 17000
$ws.Range("L11").Value = 17000
$ws.Range("N11").Value = -17284

# WVR row 113 (@@ -55153,25 +55159,25 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 970.6842
$ws.Range("I113").Value = 596.4545000000001
$ws.Range("J113").Value = 1485.25
$ws.Range("K113").Value = 1789.3635
$ws.Range("L113").Value = 4455.75
$ws.Range("M113").Value = 380.6364999999998
$ws.Range("N113").Value = -8795.75

# WVR row 132 (@@ -56069,22 +56075,22 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 12544.677
$ws.Range("I132").Value = 14858.411
$ws.Range("K132").Value = 44575.233
$ws.Range("M132").Value = -42045.233
